# Results for GraphSAGE supervised.
#
# 1) "GraphSAGE+XGBoost" sheet gains a Test_or_Val column (all "Val").
# 2) "GraphSAGE Supervised" sheet gets fully populated with the
#    gcn_small_0.0100 Loss / F1_Micro / F1_Macro results for both the
#    Test and Val splits.
# 3) A couple of cosmetic view changes: node2vec+XGBoost scrolls back to
#    the top, and the GRAND sheet becomes the active tab/selection.

$wb = $excel.ActiveWorkbook

$xlCenter = -4108

# ---------------------------------------------------------------------
# 1. "GraphSAGE Supervised" sheet: build the new results table.
# ---------------------------------------------------------------------
$wsSup = $wb.Worksheets.Item("GraphSAGE Supervised")

# Header row.
$wsSup.Cells.Item(1, 1).Value2 = "Model"
$wsSup.Cells.Item(1, 2).Value2 = "Identifier"
$wsSup.Cells.Item(1, 3).Value2 = "Loss"
$wsSup.Cells.Item(1, 4).Value2 = "F1_Micro"
$wsSup.Cells.Item(1, 5).Value2 = "F1_Macro"
$wsSup.Cells.Item(1, 6).Value2 = "Test_or_Val"

$wsSup.Range("A1:E1").HorizontalAlignment = $xlCenter

# Data rows: Test split first, then Val split, each covering the five
# model variants.
$models = @("gcn", "graphsage_maxpool", "graphsage_mean", "graphsage_meanpool", "graphsage_seq")

$testLoss   = @(0.30408000000000002, 0.23588999999999999, 0.24249999999999999, 0.21962000000000001, 0.23255000000000001)
$testMicro  = @(0.85487999999999997, 0.88771, 0.89446000000000003, 0.90105999999999997, 0.90237000000000001)
$testMacro  = @(0.48658000000000001, 0.70133999999999996, 0.75022999999999995, 0.75531000000000004, 0.77325999999999995)

$valLoss    = @(0.32694000000000001, 0.2477, 0.24138999999999999, 0.21364, 0.22484000000000001)
$valMicro   = @(0.85348999999999997, 0.88097000000000003, 0.90456999999999999, 0.91056000000000004, 0.90859999999999996)
$valMacro   = @(0.48638999999999999, 0.69376000000000004, 0.77881, 0.78483000000000003, 0.78715000000000002)

$row = 2
for ($i = 0; $i -lt $models.Length; $i++) {
    $wsSup.Cells.Item($row, 1).Value2 = $models[$i]
    $wsSup.Cells.Item($row, 2).Value2 = "gcn_small_0.0100"
    $wsSup.Cells.Item($row, 3).Value2 = $testLoss[$i]
    $wsSup.Cells.Item($row, 4).Value2 = $testMicro[$i]
    $wsSup.Cells.Item($row, 5).Value2 = $testMacro[$i]
    $wsSup.Cells.Item($row, 6).Value2 = "Test"
    $row++
}
for ($i = 0; $i -lt $models.Length; $i++) {
    $wsSup.Cells.Item($row, 1).Value2 = $models[$i]
    $wsSup.Cells.Item($row, 2).Value2 = "gcn_small_0.0100"
    $wsSup.Cells.Item($row, 3).Value2 = $valLoss[$i]
    $wsSup.Cells.Item($row, 4).Value2 = $valMicro[$i]
    $wsSup.Cells.Item($row, 5).Value2 = $valMacro[$i]
    $wsSup.Cells.Item($row, 6).Value2 = "Val"
    $row++
}

# Data rows (2:11) columns B:F are centered both ways, same as the rest
# of the workbook's result tables; column A is left at the default style.
$wsSup.Range("B2:F11").HorizontalAlignment = $xlCenter
$wsSup.Range("B2:F11").VerticalAlignment = $xlCenter

$wsSup.Columns.Item(1).AutoFit() | Out-Null
$wsSup.Columns.Item(2).AutoFit() | Out-Null
$wsSup.Columns.Item(3).AutoFit() | Out-Null
$wsSup.Columns.Item(4).AutoFit() | Out-Null
$wsSup.Columns.Item(5).AutoFit() | Out-Null
$wsSup.Columns.Item(6).AutoFit() | Out-Null

$wsSup.Range("F13").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. "GraphSAGE+XGBoost" sheet: tag every row with Test_or_Val = "Val".
# ---------------------------------------------------------------------
$wsXgb = $wb.Worksheets.Item("GraphSAGE+XGBoost")

$wsXgb.Cells.Item(1, 7).Value2 = "Test_or_Val"
for ($r = 2; $r -le 7; $r++) {
    $wsXgb.Cells.Item($r, 7).Value2 = "Val"
}
$wsXgb.Range("G2:G7").HorizontalAlignment = $xlCenter
$wsXgb.Range("G2:G7").VerticalAlignment = $xlCenter

$wsXgb.Range("G1").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. "node2vec+XGBoost" sheet: scroll view back to the top.
# ---------------------------------------------------------------------
$wsN2v = $wb.Worksheets.Item("node2vec+XGBoost")
$wsN2v.Range("A264").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. "GRAND" sheet becomes the active tab/selection.
# ---------------------------------------------------------------------
$wsGrand = $wb.Worksheets.Item("GRAND")
$wsGrand.Activate() | Out-Null
$wsGrand.Range("E13").Select() | Out-Null
